$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.407.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.646.79'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.93'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.76'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.645.70'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +8.23%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.09'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.129.89'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.319.31'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.658.95'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.36'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '364.12'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.02%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.38'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.89'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.80'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.85%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.84'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.776.40'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '574.52'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.08'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.88'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.83%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.77'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.375'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.89'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.38'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("B43").Value = 'BabyDogeCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₆0337'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.67'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.75'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.62'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '156.87'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.52%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.71'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.16%  '
